$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.184569597244263
$ws.Range("B1").Value = 2.196457624435425
$ws.Range("C1").Value = 10.52311229705811
$ws.Range("D1").Value = 2.570541143417358
$ws.Range("E1").Value = 1.238157749176025
